$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (Through 2022-04-11 -> Through 2022-04-12)
$ws.Name = "Through 2022-04-12"

# Update header cell text for the "April 2022" column
$ws.Range("B1").Value = "April 2022 (through April 12)"

# Update/insert data values per neighborhood row
# Row 2 - Austin
$ws.Range("B2").Value = 3
$ws.Range("J2").Value = 5
$ws.Range("AD2").Value = 2

# Row 3 - Englewood
$ws.Range("F3").Value = 1

# Row 5 - Garfield Park
$ws.Range("B5").Value = 4
$ws.Range("J5").Value = 3

# Row 14 - Wicker Park
$ws.Range("B14").Value = 1

# Row 16 - Washington Heights
$ws.Range("B16").Value = 2
$ws.Range("F16").Value = 1

# Row 19 - Lake View
$ws.Range("B19").Value = 2

# Row 20 - Near South Side
$ws.Range("B20").Value = 1

# Row 26 - South Shore
$ws.Range("F26").Value = 3

# Row 32 - Roseland
$ws.Range("F32").Value = 2
$ws.Range("J32").Value = 2

# Row 44 - Brighton Park
$ws.Range("J44").Value = 1

# Row 51 - Hyde Park
$ws.Range("F51").Value = 3

# Row 93 - West Pullman
$ws.Range("J93").Value = 1

# Row 94 - West Ridge
$ws.Range("F94").Value = 2
$ws.Range("R94").Value = 1
